$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the accent on "Medellin" -> "Medellín" in row 3 (Ciudad column)
$ws.Range("D3").Value = "Medellín"

# New row 6: Día Relax / Cartagena
$ws.Range("B6").NumberFormat = "@"
$ws.Range("A6").Value = "Día Relax"
$ws.Range("B6").Value = "120.0"
$ws.Range("C6").Value = "Día de Sol"
$ws.Range("D6").Value = "Cartagena"
$ws.Range("E6").Value = "Piscina, Spa, Yoga"
$ws.Range("F6").Value = "Sí"
$ws.Range("G6").Value = "Sí"

# New row 7: Aventura Extrema / Medellín
$ws.Range("B7").NumberFormat = "@"
$ws.Range("A7").Value = "Aventura Extrema"
$ws.Range("B7").Value = "150.0"
$ws.Range("C7").Value = "Día de Sol"
$ws.Range("D7").Value = "Medellín"
$ws.Range("E7").Value = "Senderismo, Rappel, Kayak"
$ws.Range("F7").Value = "No"
$ws.Range("G7").Value = "Sí"

# New row 8: Día Familiar / Bogotá
$ws.Range("B8").NumberFormat = "@"
$ws.Range("A8").Value = "Día Familiar"
$ws.Range("B8").Value = "100.0"
$ws.Range("C8").Value = "Día de Sol"
$ws.Range("D8").Value = "Bogotá"
$ws.Range("E8").Value = "Juegos Infantiles, Piscina"
$ws.Range("F8").Value = "Sí"
$ws.Range("G8").Value = "No"
